$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column ("D") cells hold plain text (e.g. "1.00", "206.83").
# Force text number-format on each one individually before writing so
# Excel does not silently reinterpret the text as a numeric value
# (union ranges such as "D2,D3" only apply formatting to the first area,
# so each cell is set one at a time).
$priceCells = @("D2", "D3", "D5", "D10", "D11", "D12", "D13", "D15", "D16", "D17", "D18", "D20", "D21", "D22", "D23", "D25", "D26", "D28", "D29", "D31", "D32", "D33", "D36", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D48")
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = '25.253.72'
$ws.Range("E2").Value = '  -2.84%  '
$ws.Range("D3").Value = '1.555.77'
$ws.Range("E3").Value = '  -4.40%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = '206.83'
$ws.Range("E6").Value = '  -0.20%  '
$ws.Range("E7").Value = '  -4.64%  '
$ws.Range("E8").Value = '  -1.55%  '
$ws.Range("E9").Value = '  -3.04%  '
$ws.Range("D10").Value = '17.81'
$ws.Range("E10").Value = '  -3.46%  '
$ws.Range("D11").Value = '0.0781'
$ws.Range("E11").Value = '  -0.94%  '
$ws.Range("D12").Value = '1.770.94'
$ws.Range("E12").Value = '  -4.46%  '
$ws.Range("D13").Value = '1.562.62'
$ws.Range("E13").Value = '  -3.87%  '
$ws.Range("E14").Value = '  -4.34%  '
$ws.Range("D15").Value = '0.506'
$ws.Range("E15").Value = '  -3.97%  '
$ws.Range("D16").Value = '25.243.79'
$ws.Range("E16").Value = '  -2.90%  '
$ws.Range("D17").Value = '58.97'
$ws.Range("E17").Value = '  -4.12%  '
$ws.Range("D18").Value = '0.0₃0708'
$ws.Range("E18").Value = '  -4.65%  '
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").Value = '185.59'
$ws.Range("E20").Value = '  -3.46%  '
$ws.Range("D21").Value = '4.12'
$ws.Range("E21").Value = '  -3.03%  '
$ws.Range("D22").Value = '9.28'
$ws.Range("E22").Value = '  -2.79%  '
$ws.Range("D23").Value = '5.85'
$ws.Range("E23").Value = '  -3.59%  '
$ws.Range("E24").Value = '  -4.02%  '
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.22%  '
$ws.Range("D26").Value = '140.80'
$ws.Range("E26").Value = '  -2.56%  '
$ws.Range("E27").Value = '  -4.75%  '
$ws.Range("D28").Value = '14.87'
$ws.Range("E28").Value = '  -2.48%  '
$ws.Range("D29").Value = '6.40'
$ws.Range("E29").Value = '  -4.75%  '
$ws.Range("E30").Value = '  -6.88%  '
$ws.Range("D31").Value = '0.0467'
$ws.Range("E31").Value = '  -3.16%  '
$ws.Range("D32").Value = '3.04'
$ws.Range("E32").Value = '  -3.13%  '
$ws.Range("D33").Value = '2.97'
$ws.Range("E33").Value = '  -4.83%  '
$ws.Range("E34").Value = '  -2.79%  '
$ws.Range("E35").Value = '  -3.80%  '
$ws.Range("D36").Value = '1.084.98'
$ws.Range("E36").Value = '  -3.45%  '
$ws.Range("E37").Value = '  -0.28%  '
$ws.Range("E38").Value = '  -3.22%  '
$ws.Range("D39").Value = '0.494'
$ws.Range("E39").Value = '  -4.82%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = '0.763'
$ws.Range("E40").Value = '  -10.24%  '
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '2.24'
$ws.Range("E41").Value = '  -7.56%  '
$ws.Range("D42").Value = '0.797'
$ws.Range("E42").Value = '  +5.87%  '
$ws.Range("D43").Value = '92.78'
$ws.Range("E43").Value = '  -5.70%  '
$ws.Range("D44").Value = '5.04'
$ws.Range("E44").Value = '  -1.05%  '
$ws.Range("D45").Value = '1.685.69'
$ws.Range("E45").Value = '  -4.41%  '
$ws.Range("E46").Value = '  -2.57%  '
$ws.Range("E47").Value = '  -2.02%  '
$ws.Range("D48").Value = '52.34'
$ws.Range("E48").Value = '  -3.91%  '
$ws.Range("E49").Value = '  -4.10%  '
$ws.Range("E50").Value = '  -0.39%  '
$ws.Range("E51").Value = '  -2.21%  '
